$d = $word.ActiveDocument

# 1. Update the main title (appears as Heading1 and again as bold text near the end)
$d.Content.Find.Execute(
    "Play Malice Free: A Dark Twist on [Alice in] Wonderland", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Play Malice Free - Dark and Intense Online Slot", 2) | Out-Null

# 2. "Adheres to slot game conventions" -> "Adheres to online slot conventions"
$d.Content.Find.Execute(
    "Adheres to slot game conventions", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Adheres to online slot conventions", 2) | Out-Null

# 3. "No progressive jackpot feature" -> "Limited number of paylines"
$d.Content.Find.Execute(
    "No progressive jackpot feature", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Limited number of paylines", 2) | Out-Null

# 4. "Limited bonus game options" -> "May not appeal to players not interested in Alice in Wonderland theme"
$d.Content.Find.Execute(
    "Limited bonus game options", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "May not appeal to players not interested in Alice in Wonderland theme", 2) | Out-Null

# 5. Update the italic summary/meta paragraph at the very end
$d.Content.Find.Execute(
    "Explore a dark twist on Alice in Wonderland in Malice, an online slot game. Play for free and experience the thrill of bonus symbols and free spins.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Read the review of Malice, an online slot game with dark and intense ambience. Play for free!", 2) | Out-Null

# 6. Remove the "Inspired by Alice in Wonderland" bullet entirely from the "What we like" list
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Inspired by Alice in Wonderland") {
        $para.Range.Delete()
        break
    }
}

# 7. Insert a new bullet "Graphics and sound design work wonderfully together" right before
#    "Special symbols enhance gameplay experience" in the "What we like" list
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Special symbols enhance gameplay experience") {
        $para.Range.InsertParagraphBefore()
        $newPara = $d.Paragraphs($i)
        $newPara.Range.InsertBefore("Graphics and sound design work wonderfully together")
        break
    }
}
